$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.114.53'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').Value = '2.495.03'
$ws.Range('E3').Value = '  -1.38%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.78'
$ws.Range('E5').Value = '  -1.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.86'
$ws.Range('E6').Value = '  -3.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.518'
$ws.Range('E7').Value = '  -1.71%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  -3.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.94'
$ws.Range('E10').Value = '  -4.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.27'
$ws.Range('E11').Value = '  -0.96%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0802'
$ws.Range('E12').Value = '  -3.15%  '

$ws.Range('E13').Value = '  +0.25%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.10'
$ws.Range('E14').Value = '  -2.89%  '

$ws.Range('D15').Value = '2.889.04'
$ws.Range('E15').Value = '  -1.35%  '

$ws.Range('D16').Value = '2.500.32'
$ws.Range('E16').Value = '  -1.12%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.829'
$ws.Range('E17').Value = '  -3.59%  '

$ws.Range('D18').Value = '48.040.82'
$ws.Range('E18').Value = '  -0.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.01'
$ws.Range('E19').Value = '  +11.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.81'
$ws.Range('E20').Value = '  -3.75%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.58'

$ws.Range('D22').Value = '0.0₃0929'
$ws.Range('E22').Value = '  -2.52%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.07'
$ws.Range('E23').Value = '  -1.71%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.74'
$ws.Range('E24').Value = '  -0.98%  '

$ws.Range('E25').Value = '  -2.60%  '

$ws.Range('E26').Value = '  +0.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.74'
$ws.Range('E27').Value = '  -1.89%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -2.05%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.74'
$ws.Range('E29').Value = '  -4.13%  '

$ws.Range('E30').Value = '  -2.25%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.61'
$ws.Range('E31').Value = '  -3.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.30'
$ws.Range('E32').Value = '  -0.88%  '

$ws.Range('E33').Value = '  +0.00%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.08'
$ws.Range('E34').Value = '  -4.20%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.28'
$ws.Range('E35').Value = '  -2.56%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0772'
$ws.Range('E36').Value = '  -2.90%  '

$ws.Range('E37').Value = '  -2.39%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.59'
$ws.Range('E38').Value = '  -3.65%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.88'
$ws.Range('E39').Value = '  -4.43%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '122.86'
$ws.Range('E40').Value = '  +2.68%  '

$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.44'
$ws.Range('E41').Value = '  +0.92%  '

$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.110'
$ws.Range('E42').Value = '  -2.01%  '

$ws.Range('E43').Value = '  +1.27%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0303'
$ws.Range('E44').Value = '  +0.28%  '

$ws.Range('D45').Value = '2.000.90'
$ws.Range('E45').Value = '  -0.75%  '

$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.91'
$ws.Range('E47').Value = '  +0.57%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.94'
$ws.Range('E49').Value = '  -2.51%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.20'
$ws.Range('E50').Value = '  -1.42%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.68'
$ws.Range('E51').Value = '  -1.43%  '
